$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row above the old "closing" row of the table (row 19),
#    shifting everything below (blank spacer rows + the two signature rows) down by one.
$ws.Rows("19:19").Insert()

# 2) The row that used to be the last data row (row 18, period 2508) gets copied
#    down into the freshly inserted row 19 -- this carries over both its values
#    (CC / 19772276 / LUIS ALBERTO AREVALO ARENILLA / 2508 / 128116 / 3202875)
#    and its "closing" (bottom-border) formatting.
$ws.Range("B18:J18").Copy($ws.Range("B19:J19"))

# 3) Row 18 becomes a regular interior row of the table now, so it picks up the
#    same formatting as row 17 (format-only paste -- the values in row 18 stay put).
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4) New row 19 represents the new period (2509) -- same worker, same salary base.
$ws.Range("E19").Value = "2509"
$ws.Range("F19").Value = 128116
$ws.Range("G19").Value = 3202875

# 5) Center the "Periodo Mora" column for all the data rows.
$ws.Range("E16:E19").HorizontalAlignment = -4108

# 6) Refresh the summary figures at the top of the account statement.
$ws.Range("E11").Value = 469758
$ws.Range("F13").Value = 4
